$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C7 was blank; fill it with the literal text "=" (not a formula).
# Assigning "=" directly would be parsed as a formula, and prefixing with
# an apostrophe stamps a quotePrefix style onto the cell, which would
# change its style index. Round-trip the literal value through a scratch
# cell + PasteSpecial(values-only) so the text lands without disturbing
# C7's existing style.
$scratch = $ws.Range("Z1")
$scratch.Value = "'="
$scratch.Copy()
$ws.Range("C7").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

# D13 was blank; fill it with the literal text "String" (matches C13/E13).
$ws.Range("D13").Value = "String"

# Restore the original selection/active cell.
$ws.Range("D16").Select() | Out-Null
